$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.121.43"
$ws.Range("E2").Value = "  -0.95%  "
$ws.Range("D3").Value = "1.894.45"
$ws.Range("E3").Value = "  -0.62%  "
$ws.Range("E4").Value = "  -0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.90"
$ws.Range("E6").Value = "  -1.28%  "
$ws.Range("E7").Value = "  -0.46%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.36"
$ws.Range("E8").Value = "  +1.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.352"
$ws.Range("E9").Value = "  -0.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.07"
$ws.Range("E10").Value = "  -1.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0747"
$ws.Range("E11").Value = "  +2.90%  "
$ws.Range("E12").Value = "  -1.35%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "12.98"
$ws.Range("E13").Value = "  +3.38%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.168.43"
$ws.Range("E14").Value = "  -0.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.725"
$ws.Range("E15").Value = "  +1.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.95"
$ws.Range("E16").Value = "  +1.21%  "
$ws.Range("D17").Value = "1.910.23"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").Value = "35.109.37"
$ws.Range("E18").Value = "  -0.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.71"
$ws.Range("E19").Value = "  +0.88%  "
$ws.Range("D20").Value = "0.0₃0828"
$ws.Range("E20").Value = "  +0.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "251.04"
$ws.Range("E21").Value = "  +3.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.94"
$ws.Range("E22").Value = "  +0.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.00"
$ws.Range("E23").Value = "  -1.01%  "
$ws.Range("E24").Value = "  -0.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.42"
$ws.Range("E25").Value = "  +4.36%  "
$ws.Range("E26").Value = "  -3.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "167.74"
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.53"
$ws.Range("E28").Value = "  -1.85%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.36"
$ws.Range("E29").Value = "  -2.69%  "
$ws.Range("E30").Value = "  -2.95%  "
$ws.Range("D31").Value = "4.128.37"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.30"
$ws.Range("E32").Value = "  +1.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0588"
$ws.Range("E33").Value = "  +2.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.94"
$ws.Range("E34").Value = "  +3.74%  "
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.23"
$ws.Range("E35").Value = "  +0.58%  "
$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.56"
$ws.Range("E36").Value = "  +6.67%  "
$ws.Range("E37").Value = "  -0.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.843"
$ws.Range("E38").Value = "  -8.73%  "
$ws.Range("E39").Value = "  -0.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.50"
$ws.Range("E40").Value = "  +5.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.66"
$ws.Range("E41").Value = "  +2.63%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0663"
$ws.Range("E42").Value = "  +2.11%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0213"
$ws.Range("E43").Value = "  +1.80%  "
$ws.Range("E44").Value = "  -2.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.39"
$ws.Range("E45").Value = "  -0.85%  "
$ws.Range("D46").Value = "1.297.61"
$ws.Range("E46").Value = "  -4.12%  "
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.74"
$ws.Range("E48").Value = "  -1.58%  "
$ws.Range("E49").Value = "  +8.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.52"
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "12.02"
$ws.Range("E51").Value = "  -1.26%  "
